# Update "want-to-go" counts (column F) on the "展览" (exhibitions),
# "演出" (performances) and "全部类型" (all types) sheets to reflect
# the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 120
$wsExhibit.Range("F3").Value = 2138
$wsExhibit.Range("F4").Value = 23
$wsExhibit.Range("F5").Value = 11131
$wsExhibit.Range("F8").Value = 304
$wsExhibit.Range("F10").Value = 11033
$wsExhibit.Range("F13").Value = 39
$wsExhibit.Range("F14").Value = 1717
$wsExhibit.Range("F15").Value = 5537
$wsExhibit.Range("F17").Value = 3431

# --- 演出 (Performances) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 566

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 120
$wsAll.Range("F3").Value = 2138
$wsAll.Range("F4").Value = 566
$wsAll.Range("F5").Value = 23
$wsAll.Range("F7").Value = 11131
$wsAll.Range("F10").Value = 304
$wsAll.Range("F12").Value = 11033
$wsAll.Range("F15").Value = 39
$wsAll.Range("F16").Value = 1717
$wsAll.Range("F17").Value = 5537
$wsAll.Range("F19").Value = 3431
